# approvalProcess.xlsx edit:
#  - Update the item-description test-data cell on the FinanceApprover
#    sheet from "REPOFLOR 100 MG" to "DESKTOP".
#  - Move the sheet selection/scroll position from U1 (with J1 pinned as
#    the top-left visible cell) to E10 (back at the sheet's natural
#    top-left), matching the updated "Test steps" in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FinanceApprover")

# New test data value for the item description column (row 2).
$ws.Range("C2").Value = "DESKTOP"

# Re-anchor the view: select E10, which also resets the scrolled
# top-left cell back to the sheet default.
$ws.Activate()
$ws.Range("E10").Select()
